# Generate Report for Archive
#
# The localization-status report re-sorts the file list: the row for
# 949f0782-7007-4ea7-b6e0-d9a395d6c1cb.md (previously last, row 7) moves up
# to row 4 (right after 70fbc62f...), and the rows in between
# (d0826aa6, 0460c0ac, 474c97d6) each shift down by one row, carrying all of
# their own column data (status / handoff file / handoff datetime / etc.)
# and their cell hyperlinks with them.
#
# This happens identically on all three sheets: "Overview" (cols A-C),
# "zh-cn" (cols A-D, G-H, hyperlinks on A & C), "de-de" (same as zh-cn).

$wb = $excel.ActiveWorkbook

# Rows 4..7 or a given sheet are read first (old order), then written back
# in the new order: new[4]=old[7], new[5]=old[4], new[6]=old[5], new[7]=old[6]
$oldOrder = @(4, 5, 6, 7)
$newOrder = @(7, 4, 5, 6)   # newOrder[i] = source row for destination row $oldOrder[i]

function Move-ReportRows {
    param(
        $ws,
        [int[]]$cols          # 1-based column indices that hold data on these rows
    )

    # --- capture old cell values for rows 4..7, for every tracked column ---
    $old = @{}
    foreach ($r in $oldOrder) {
        $old[$r] = @{}
        foreach ($c in $cols) {
            $old[$r][$c] = $ws.Cells.Item($r, $c).Value()
        }
    }

    # --- capture old hyperlink display text, keyed by "row:col" ---
    $oldLinks = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()        # like "$A$4"
        $addr = $addr -replace '\$', ''    # -> "A4"
        if ($addr -match '^([A-Z]+)([0-9]+)$') {
            $colLetters = $Matches[1]
            $rowNum = [int]$Matches[2]
            if ($oldOrder -contains $rowNum) {
                $oldLinks["${rowNum}:${colLetters}"] = $hl.TextToDisplay
            }
        }
    }

    # --- write new values: destination row <- source row ---
    for ($i = 0; $i -lt $oldOrder.Count; $i++) {
        $destRow = $oldOrder[$i]
        $srcRow = $newOrder[$i]
        foreach ($c in $cols) {
            $ws.Cells.Item($destRow, $c).Value = $old[$srcRow][$c]
        }
    }

    # --- fix up hyperlink display text to match the row it now lives on ---
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        $addr = $addr -replace '\$', ''
        if ($addr -match '^([A-Z]+)([0-9]+)$') {
            $colLetters = $Matches[1]
            $rowNum = [int]$Matches[2]
            $idx = [array]::IndexOf($oldOrder, $rowNum)
            if ($idx -ge 0) {
                $srcRow = $newOrder[$idx]
                $key = "${srcRow}:${colLetters}"
                if ($oldLinks.ContainsKey($key)) {
                    $hl.TextToDisplay = $oldLinks[$key]
                }
            }
        }
    }
}

# Sheet "Overview": File Name (A), zh-cn (B), de-de (C) - hyperlinks on col A only
$wsOverview = $wb.Worksheets.Item("Overview")
Move-ReportRows $wsOverview @(1, 2, 3)

# Sheet "zh-cn": columns A,B,C,D,G,H populated on rows 4-7, hyperlinks on A & C
$wsZh = $wb.Worksheets.Item("zh-cn")
Move-ReportRows $wsZh @(1, 2, 3, 4, 7, 8)

# Sheet "de-de": same layout as zh-cn
$wsDe = $wb.Worksheets.Item("de-de")
Move-ReportRows $wsDe @(1, 2, 3, 4, 7, 8)
